$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(188, "Nissan Skyline GT-R (R32)", "Color Shifters"),
    @(189, "Gordon Murray Automotive T.33", "Mainlines"),
    @(190, "Porsche 911 Turbo Cabriolet", "Mainlines"),
    @(191, "89 Porsche 944 Turbo", "Mainlines"),
    @(192, "2013 SRT Viper", "Mainlines"),
    @(193, "92 BMW M3", "Mainlines")
)

$startRow = 189
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}
